$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.04063791522011
$ws.Range("D2").Value = 1.042651887187782
$ws.Range("E2").Value = 1.044294540807085
$ws.Range("F2").Value = 1.049810112539466
$ws.Range("I2").Value = 1.042997334198623
$ws.Range("J2").Value = 1.045723442190164
$ws.Range("K2").Value = 1.045427971016013
$ws.Range("L2").Value = 1.047065994778422
$ws.Range("M2").Value = 1.052566135684847
$ws.Range("N2").Value = 1.019019783251242

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.041786771687153
$ws.Range("D3").Value = 1.043527385752907
$ws.Range("E3").Value = 1.045392491129757
$ws.Range("F3").Value = 1.051181456834945
$ws.Range("I3").Value = 1.043395302136456
$ws.Range("J3").Value = 1.046517003751703
$ws.Range("K3").Value = 1.046114175405057
$ws.Range("L3").Value = 1.047974405651792
$ws.Range("M3").Value = 1.053748358689667
$ws.Range("N3").Value = 1.019295199737327

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.042529679970808
$ws.Range("D4").Value = 1.044093427959927
$ws.Range("E4").Value = 1.046102809424225
$ws.Range("F4").Value = 1.052068690613659
$ws.Range("I4").Value = 1.043651254652377
$ws.Range("J4").Value = 1.047029479096298
$ws.Range("K4").Value = 1.046557096049042
$ws.Range("L4").Value = 1.048561490795031
$ws.Range("M4").Value = 1.054512687079787
$ws.Range("N4").Value = 1.019472791575708

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.042841886002412
$ws.Range("D5").Value = 1.044331281872527
$ws.Range("E5").Value = 1.046401397755927
$ws.Range("F5").Value = 1.052441658577999
$ws.Range("I5").Value = 1.04375848450159
$ws.Range("J5").Value = 1.047244682661075
$ws.Range("K5").Value = 1.046743037619681
$ws.Range("L5").Value = 1.048808130924406
$ws.Range("M5").Value = 1.05483385800038
$ws.Range("N5").Value = 1.01954730279539

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.042894300222977
$ws.Range("D6").Value = 1.044371212151068
$ws.Range("E6").Value = 1.04645153036261
$ws.Range("F6").Value = 1.052504280156317
$ws.Range("I6").Value = 1.043776467037336
$ws.Range("J6").Value = 1.047280802185388
$ws.Range("K6").Value = 1.0467742426756
$ws.Range("L6").Value = 1.048849532928719
$ws.Range("M6").Value = 1.054887775130976
$ws.Range("N6").Value = 1.019559804863847

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.042533852123083
$ws.Range("D7").Value = 1.044096606607652
$ws.Range("E7").Value = 1.046106799290211
$ws.Range("F7").Value = 1.052073674324723
$ws.Range("I7").Value = 1.043652688925988
$ws.Range("J7").Value = 1.047032355602588
$ws.Range("K7").Value = 1.046559581638763
$ws.Range("L7").Value = 1.048564787080712
$ws.Range("M7").Value = 1.054516979178214
$ws.Range("N7").Value = 1.019473787780977

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.041026276949729
$ws.Range("D8").Value = 1.042947862582932
$ws.Range("E8").Value = 1.044665625809464
$ws.Range("F8").Value = 1.050273590971796
$ws.Range("I8").Value = 1.043132152636549
$ws.Range("J8").Value = 1.045991840320402
$ws.Range("K8").Value = 1.04566010522986
$ws.Range("L8").Value = 1.047373145814164
$ws.Range("M8").Value = 1.05296580943402
$ws.Range("N8").Value = 1.019112990308743

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.03836598161752
$ws.Range("D9").Value = 1.040920038119246
$ws.Range("E9").Value = 1.042125042994611
$ws.Range("F9").Value = 1.047100575786047
$ws.Range("I9").Value = 1.042202925753397
$ws.Range("J9").Value = 1.044150513465263
$ws.Range("K9").Value = 1.044066649191878
$ws.Range("L9").Value = 1.045267768705913
$ws.Range("M9").Value = 1.050227362248472
$ws.Range("N9").Value = 1.01847244940516

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.03658979951907
$ws.Range("D10").Value = 1.039565685508056
$ws.Range("E10").Value = 1.040430523536438
$ws.Range("F10").Value = 1.044984340633691
$ws.Range("I10").Value = 1.041575342467073
$ws.Range("J10").Value = 1.042917635034224
$ws.Range("K10").Value = 1.042998592496262
$ws.Range("L10").Value = 1.043860366701022
$ws.Range("M10").Value = 1.048398143368051
$ws.Range("N10").Value = 1.018042195851316

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.035820035065357
$ws.Range("D11").Value = 1.038978636942413
$ws.Range("E11").Value = 1.039696567668038
$ws.Range("F11").Value = 1.044067732256176
$ws.Range("I11").Value = 1.041301660134387
$ws.Range("J11").Value = 1.042382503291301
$ws.Range("K11").Value = 1.042534733205166
$ws.Range("L11").Value = 1.043250021310869
$ws.Range("M11").Value = 1.047605182816228
$ws.Range("N11").Value = 1.017855121442855

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.035534007634964
$ws.Range("D12").Value = 1.038760488846125
$ws.Range("E12").Value = 1.039423909148875
$ws.Range("F12").Value = 1.043727219489562
$ws.Range("I12").Value = 1.041199710660387
$ws.Range("J12").Value = 1.042183536478319
$ws.Range("K12").Value = 1.042362225988933
$ws.Range("L12").Value = 1.04302317021444
$ws.Range("M12").Value = 1.047310503903244
$ws.Range("N12").Value = 1.0177855172843

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.035595366176127
$ws.Range("D13").Value = 1.03880728654186
$ws.Range("E13").Value = 1.039482396924717
$ws.Range("F13").Value = 1.043800262669004
$ws.Range("I13").Value = 1.041221592383633
$ws.Range("J13").Value = 1.04222622438862
$ws.Range("K13").Value = 1.042399238861699
$ws.Range("L13").Value = 1.043071836982757
$ws.Range("M13").Value = 1.047373719831483
$ws.Range("N13").Value = 1.01780045288678

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.035796394057578
$ws.Range("D14").Value = 1.038960606639123
$ws.Range("E14").Value = 1.039674030347335
$ws.Range("F14").Value = 1.044039586276474
$ws.Range("I14").Value = 1.041293238914917
$ws.Range("J14").Value = 1.042366060621215
$ws.Range("K14").Value = 1.042520477980688
$ws.Range("L14").Value = 1.043231272628181
$ws.Range("M14").Value = 1.047580827401027
$ws.Range("N14").Value = 1.01784937031836

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.035920240291676
$ws.Range("D15").Value = 1.039055059967309
$ws.Range("E15").Value = 1.039792097377467
$ws.Range("F15").Value = 1.04418703555646
$ws.Range("I15").Value = 1.041337344027991
$ws.Range("J15").Value = 1.04245219244251
$ws.Range("K15").Value = 1.042595149639197
$ws.Range("L15").Value = 1.043329487369373
$ws.Range("M15").Value = 1.047708414816071
$ws.Range("N15").Value = 1.017879494525859

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.036640871999814
$ws.Range("D16").Value = 1.039604633137257
$ws.Range("E16").Value = 1.040479229019034
$ws.Range("F16").Value = 1.045045167009821
$ws.Range("I16").Value = 1.041593465015284
$ws.Range("J16").Value = 1.042953122703017
$ws.Range("K16").Value = 1.043029348042841
$ws.Range("L16").Value = 1.043900853548088
$ws.Range("M16").Value = 1.048450750371114
$ws.Range("N16").Value = 1.018054595065447

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.037092725125418
$ws.Range("D17").Value = 1.039949203037414
$ws.Range("E17").Value = 1.040910188895203
$ws.Range("F17").Value = 1.045583376646124
$ws.Range("I17").Value = 1.041753604274731
$ws.Range("J17").Value = 1.043266997217458
$ws.Range("K17").Value = 1.043301337833944
$ws.Range("L17").Value = 1.044259005948855
$ws.Range("M17").Value = 1.048916155426902
$ws.Range("N17").Value = 1.018164224068805

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.037356219335366
$ws.Range("D18").Value = 1.040150126545109
$ws.Range("E18").Value = 1.04116153961784
$ws.Range("F18").Value = 1.045897279950945
$ws.Range("I18").Value = 1.041846824200055
$ws.Range("J18").Value = 1.043449950822688
$ws.Range("K18").Value = 1.043459851461105
$ws.Range("L18").Value = 1.044467820472281
$ws.Range("M18").Value = 1.049187532027268
$ws.Range("N18").Value = 1.018228094375889

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.037446053262548
$ws.Range("D19").Value = 1.040218626420096
$ws.Range("E19").Value = 1.041247240234934
$ws.Range("F19").Value = 1.046004308696599
$ws.Range("I19").Value = 1.041878578166548
$ws.Range("J19").Value = 1.043512312277001
$ws.Range("K19").Value = 1.043513877922488
$ws.Range("L19").Value = 1.044539005665312
$ws.Range("M19").Value = 1.049280049899295
$ws.Range("N19").Value = 1.018249859887866

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.037044252207186
$ws.Range("D20").Value = 1.039912239985531
$ws.Range("E20").Value = 1.040863953152511
$ws.Range("F20").Value = 1.045525634486914
$ws.Range("I20").Value = 1.041736442155343
$ws.Range("J20").Value = 1.043233334299385
$ws.Range("K20").Value = 1.043272169712773
$ws.Range("L20").Value = 1.044220588877653
$ws.Range("M20").Value = 1.048866230818541
$ws.Range("N20").Value = 1.018152469616626

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.035737199211203
$ws.Range("D21").Value = 1.038915460244816
$ws.Range("E21").Value = 1.039617600052091
$ws.Range("F21").Value = 1.043969112673709
$ws.Range("I21").Value = 1.041272148856875
$ws.Range("J21").Value = 1.042324887733955
$ws.Range("K21").Value = 1.042484781868283
$ws.Range("L21").Value = 1.043184326728686
$ws.Range("M21").Value = 1.047519843221733
$ws.Range("N21").Value = 1.017834968570598

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.03491480798926
$ws.Range("D22").Value = 1.038288211660939
$ws.Range("E22").Value = 1.03883376707647
$ws.Range("F22").Value = 1.042990211441491
$ws.Range("I22").Value = 1.040978541215574
$ws.Range("J22").Value = 1.041752581954865
$ws.Range("K22").Value = 1.041988508717161
$ws.Range("L22").Value = 1.042531967204916
$ws.Range("M22").Value = 1.046672515483853
$ws.Range("N22").Value = 1.017634669422253

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.035350830416913
$ws.Range("D23").Value = 1.038620778996633
$ws.Range("E23").Value = 1.03924931149128
$ws.Range("F23").Value = 1.043509170889892
$ws.Range("I23").Value = 1.041134348525566
$ws.Range("J23").Value = 1.042056079684879
$ws.Range("K23").Value = 1.042251707757078
$ws.Range("L23").Value = 1.042877873685893
$ws.Range("M23").Value = 1.047121776937543
$ws.Range("N23").Value = 1.017740915792184

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.03706615522976
$ws.Range("D24").Value = 1.039928942177337
$ws.Range("E24").Value = 1.040884845158441
$ws.Range("F24").Value = 1.045551725760333
$ws.Range("I24").Value = 1.04174419755462
$ws.Range("J24").Value = 1.043248545506017
$ws.Range("K24").Value = 1.043285349943054
$ws.Range("L24").Value = 1.044237948175199
$ws.Range("M24").Value = 1.048888789865895
$ws.Range("N24").Value = 1.01815778117746

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.039054189802779
$ws.Range("D25").Value = 1.041444710548177
$ws.Range("E25").Value = 1.042781978436634
$ws.Range("F25").Value = 1.047921019216391
$ws.Range("I25").Value = 1.042444577163544
$ws.Range("J25").Value = 1.044627473744024
$ws.Range("K25").Value = 1.044479605142367
$ws.Range("L25").Value = 1.045812725808853
$ws.Range("M25").Value = 1.05093593655524
$ws.Range("N25").Value = 1.018638611805088
